# Update handback report timestamps (regenerated report values)
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-13 16:46:25"
$wsZhCn.Range("H4").Value = "2016-03-13 16:46:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-13 16:46:29"
$wsDeDe.Range("H4").Value = "2016-03-13 16:46:54"
